# Added notes and user stories for Print 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (new story, Sprint 3) ---
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = "As a system admin, I would like to develop the back end in Python and link it with JavaScript"

# --- Row 5 (new story, number/sprint/state left blank) ---
$ws.Cells.Item(5, 3).Value = "As a user, I would like to received context spell checking"

# State for row 4, entered after row 5's story so the shared-string
# table order matches the source workbook.
$ws.Cells.Item(4, 4).Value = "In Progress"

# --- Row 3 (existing row 3 already had Story Number/Sprint filled in) ---
$ws.Cells.Item(3, 3).Value = "As a user, I would like to receive spellchecking wihtout having to click a submit button"
$ws.Cells.Item(3, 4).Value = "In Progress"

# Column C was widened (and no longer "best fit") to accommodate the
# longer story text added above.
$ws.Columns.Item(3).ColumnWidth = 73.83

# Leave the selection on the last-edited cell, like the source workbook.
$ws.Range("C5").Select()
